# Insert a new weekly price-report row into the daily "Acelga" sheet.
# The new row is inserted at row 110 (right after the row for date
# serial 44252 / row 109), pushing the existing rows 110..204 down to
# 111..205. The new row carries the same market/commune/category
# metadata as the (old) row 110, but a new reporting date (44589) and
# a new "Volumen" value (60); the rest of the price fields match the
# unaffected neighbouring data exactly as captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 110:204 down to 111:205 by inserting a blank row at 110.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with its data.
$ws.Range("A110").Value = 7
$ws.Range("B110").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C110").Value = "Ñuble"
$ws.Range("D110").Value = 44589
$ws.Range("E110").Value = 16
$ws.Range("F110").Value = 100112009
$ws.Range("G110").Value = "Acelga"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 60
$ws.Range("K110").Value = 350
$ws.Range("L110").Value = 400
$ws.Range("M110").Value = 375
$ws.Range("N110").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O110").Value = "Provincia de Diguillín"
$ws.Range("P110").Value = 375
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"
